$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F1 previously held a HYPERLINK() formula whose displayed text was the old
# page title. It is now a plain text cell carrying the new metadata-portal
# label (the existing hyperlink relationship on F1 is left untouched).
$ws.Range("F1").Value = "Metadata - Single European Sky Portal"

# The support-contact address changed (and lost its former ALL-CAPS/NSA prefix).
$ws.Range("F2").Value = "pru-support@eurocontrol.int"

# The hyperlink-style font used by F1 switched from Calibri/#396EA2 to Arial/#0000FF.
$ws.Range("F1").Font.Name = "Arial"
$ws.Range("F1").Font.Color = 16711680
